# -------------------------------------------------------------------------
# Adds a "Summary" sheet that aggregates the per-fold accuracy already
# computed on "Base Fold" / "Experiment Fold", and wires each of those two
# sheets up with a helper "Fold Acc" column (K) that pulls each fold's
# average accuracy so the new summary table can reference it directly.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$baseFold = $wb.Worksheets.Item("Base Fold")
$expFold  = $wb.Worksheets.Item("Experiment Fold")

# --- Base Fold: add column K "Fold Acc" ----------------------------------
$baseFold.Range("K2").Value = "Fold Acc"

$baseFold.Range("K3").Formula  = "=H13"
$baseFold.Range("K4").Formula  = "=H25"
$baseFold.Range("K5").Formula  = "=H38"
$baseFold.Range("K6").Formula  = "=H50"
$baseFold.Range("K7").Formula  = "=H62"
$baseFold.Range("K8").Formula  = "=H74"
$baseFold.Range("K9").Formula  = "=H86"
$baseFold.Range("K10").Formula = "=H98"
$baseFold.Range("K11").Formula = "=H110"
$baseFold.Range("K12").Formula = "=H122"

$baseFold.Range("L13").ClearContents()
$baseFold.Range("J13").Value = "Average"
$baseFold.Range("K13").Formula = "=AVERAGE(K3:K12)"
$baseFold.Range("J14").Value = "Stdev"
$baseFold.Range("K14").Formula = "=STDEV.P(K3:K12)"

$baseFold.Activate()
$excel.ActiveWindow.Zoom = 115
$baseFold.Range("K2:K12").Select()

# --- Experiment Fold: add column K "Fold Acc" -----------------------------
$expFold.Range("K2").Value = "Fold Acc"

$expFold.Range("K3").Formula  = "=H13"
$expFold.Range("K4").Formula  = "=H25"
$expFold.Range("K5").Formula  = "=H37"
$expFold.Range("K6").Formula  = "=H49"
$expFold.Range("K7").Formula  = "=H61"
$expFold.Range("K8").Formula  = "=H73"
$expFold.Range("K9").Formula  = "=H85"
$expFold.Range("K10").Formula = "=H97"
$expFold.Range("K11").Formula = "=H109"
$expFold.Range("K12").Formula = "=H121"

$expFold.Range("L13").ClearContents()
$expFold.Range("J13").Value = "Average"
$expFold.Range("K13").Formula = "=AVERAGE(K3:K12)"
$expFold.Range("J14").Value = "Stdev"
$expFold.Range("K14").Formula = "=STDEV.P(K3:K12)"

$expFold.Activate()
$excel.ActiveWindow.Zoom = 115
$expFold.Range("K2:K12").Select()

# --- New "Summary" sheet, placed after "Experiment Fold" -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$summary = $wb.Worksheets.Add($null, $lastSheet)
$summary.Name = "Summary"

$summary.Range("A2").Value = "Fold"
$summary.Range("C2").Value = "With Subgraph`n(concat)"
$summary.Range("B2").Value = "Base (GCN)"

$headerRange = $summary.Range("A2:C2")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$summary.Range("C2").WrapText = $true

for ($i = 0; $i -lt 10; $i++) {
    $row = 3 + $i
    $summary.Cells.Item($row, 1).Value = $i + 1
    $summary.Cells.Item($row, 2).Formula = "=Base Fold!K" + (3 + $i)
    $summary.Cells.Item($row, 3).Formula = "=Experiment Fold!K" + (3 + $i)
}

$summary.Range("A3:A12").Borders.LineStyle = 1

$dataRange = $summary.Range("B3:C12")
$dataRange.NumberFormat = "0.000"
$dataRange.Borders.LineStyle = 1

$summary.Range("A13").Value = "Average"
$summary.Range("B13").Formula = "=AVERAGE(B3:B12)"
$summary.Range("C13").Formula = "=AVERAGE(C3:C12)"

$summary.Range("A14").Value = "Stdev"
$summary.Range("B14").Formula = "=STDEV.P(B3:B12)"
$summary.Range("C14").Formula = "=STDEV.P(C3:C12)"

$summary.Range("A13:A14").Font.Bold = $true
$summary.Range("A13:A14").HorizontalAlignment = -4152
$summary.Range("A13:A14").Borders.LineStyle = 1

$summary.Range("B13:C14").NumberFormat = "0.000"
$summary.Range("B13:C14").Borders.LineStyle = 1
$summary.Range("C13").Font.Bold = $true

$summary.Columns.Item(2).ColumnWidth = 13.14
$summary.Columns.Item(3).ColumnWidth = 21.3

$summary.Rows.Item(2).RowHeight = 30

$summary.PageSetup.Orientation = 1

$summary.Range("E15").Select()

Write-Output "Summary sheet created"
